$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '79.638.83'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +4.28%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.204.77'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +5.22%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.09%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '210.54'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +6.13%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '635.41'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.74%  '

$ws.Range('E7').Value = '  -0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.246'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +19.61%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.601'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +9.63%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.199.89'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.12%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.600'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +36.87%  '

$ws.Range('E12').Value = '  +35.62%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.166'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.27%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.41'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.22%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '3.786.34'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +4.98%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '32.14'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +11.52%  '

$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '79.524.09'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.17%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.197.56'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +5.13%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.53'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +7.39%  '

$ws.Range('B20').Value = 'SuiNetwork'
$ws.Range('C20').Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.01'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +25.64%  '

$ws.Range('B21').Value = 'Uniswap'
$ws.Range('C21').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '9.35'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +3.96%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '440.26'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +15.53%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.22'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +19.42%  '

$ws.Range('E24').Value = '  +10.92%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '3.360.46'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +5.38%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '77.22'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +6.45%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.89'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +11.41%  '

$ws.Range('E28').Value = '  -0.13%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '0.0000125'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +15.44%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.16'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +10.63%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.998'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.12%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.53'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.02%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '557.00'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +12.36%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.158'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +34.04%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +6.37%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '23.08'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +12.01%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.122'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +16.59%  '

$ws.Range('E38').Value = '  -0.09%  '

$ws.Range('E39').Value = '  +8.83%  '

$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '163.39'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.03%  '

$ws.Range('E41').Value = '  -0.01%  '

$ws.Range('B42').Value = 'RenderToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.64'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +10.54%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '192.16'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.03%  '

$ws.Range('B44').Value = 'USDe'
$ws.Range('C44').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.00'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.05%  '

$ws.Range('E45').Value = '  +11.68%  '

$ws.Range('B46').Value = 'dogwifhat'
$ws.Range('C46').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.70'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +11.91%  '

$ws.Range('B47').Value = 'Mantle'
$ws.Range('C47').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.800'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.21%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.35'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +7.47%  '

$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '43.32'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +3.42%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.643'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +7.35%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '25.59'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +15.31%  '

